$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Degree column (C) cells: replace long-form degree text with abbreviations
# "Minor Misdemeanor" -> "MM", "Unclassified Misdemeanor" -> "UCM"
$ws.Range("C2").Value = "MM"
$ws.Range("C3").Value = "MM"
$ws.Range("C7").Value = "MM"
$ws.Range("C8").Value = "MM"
$ws.Range("C9").Value = "MM"
$ws.Range("C10").Value = "UCM"
$ws.Range("C12").Value = "UCM"
$ws.Range("C13").Value = "MM"
$ws.Range("C14").Value = "UCM"
$ws.Range("C15").Value = "UCM"
$ws.Range("C16").Value = "MM"
$ws.Range("C17").Value = "MM"
$ws.Range("C18").Value = "MM"
$ws.Range("C19").Value = "MM"
$ws.Range("C20").Value = "MM"
$ws.Range("C21").Value = "UCM"
$ws.Range("C23").Value = "MM"
$ws.Range("C24").Value = "MM"
$ws.Range("C25").Value = "UCM"
$ws.Range("C26").Value = "UCM"
$ws.Range("C27").Value = "MM"
$ws.Range("C28").Value = "MM"
$ws.Range("C29").Value = "MM"
$ws.Range("C30").Value = "MM"
$ws.Range("C31").Value = "MM"
$ws.Range("C32").Value = "MM"
$ws.Range("C33").Value = "MM"
$ws.Range("C34").Value = "MM"
$ws.Range("C35").Value = "MM"
$ws.Range("C36").Value = "MM"
$ws.Range("C37").Value = "MM"
$ws.Range("C38").Value = "UCM"

# Restore the active cell selection used when the workbook was last saved
$ws.Range("C39").Select() | Out-Null
